# Fixed workflow: the first 4 "Cutoff" rows (Cutoff = 1..4) were dropped
# from each reaction-sensitivity table. The remaining rows (formerly
# Cutoff = 5..19) shift up so they directly follow the header, while the
# left-most index column (A, a plain 0-based row counter) keeps counting
# 0..14 - it does not carry over the old values. Each sheet shrinks from
# 19 data rows (A1:C20) to 15 data rows (A1:C16).

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {

    # Snapshot the current Cutoff (col B) / Reaction_number (col C) pairs
    # for data rows 6..20 (i.e. Cutoff 5..19) before touching anything.
    $kept = New-Object System.Collections.ArrayList
    for ($r = 6; $r -le 20; $r++) {
        $b = $ws.Cells.Item($r, 2).Value2
        $c = $ws.Cells.Item($r, 3).Value2
        [void]$kept.Add(@($b, $c))
    }

    # Write the kept pairs into rows 2..16, directly under the header;
    # column A already holds 0..14 there and is left untouched.
    for ($i = 0; $i -lt $kept.Count; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $kept[$i][0]
        $ws.Cells.Item($row, 3).Value = $kept[$i][1]
    }

    # Drop the now-duplicated trailing rows 17..20, shrinking the sheet
    # back down to A1:C16.
    $ws.Range("A17:C20").EntireRow.Delete()
}
